$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Resolving-Mac" sending-cluster rows (old rows 14-17);
# the dataset was recomputed with new TPM values and that sending cluster
# is no longer present.
$ws.Rows("14:17").Delete()

# Updated values (new TPM) for the remaining data rows (2-13).
$newValues = @(
    @("ECs", "Mmrn2", "Cd93", "ECs", 3, 1, 53.20407366666667, 159.612221, 0.9416818042520982, 0.9416818042520981, 3, 1, 281.0920463333333, 843.2761389999999, 0.8291026083535286, 0.8291026083535286, 14955.24194023275, 134597.1774620947, 0.7807508401444716, 0.7807508401444715),
    @("ECs", "Mmrn2", "Cd93", "FAPs", 3, 1, 53.20407366666667, 159.612221, 0.9416818042520982, 0.9416818042520981, 2, 0.6666666666666666, 0.449122, 1.347366, 0.001324719879221983, 0.001324719879221983, 23.89511997332067, 215.056079759886, 0.001247464605994379, 0.001247464605994379),
    @("ECs", "Mmrn2", "Cd93", "MuSCs", 3, 1, 53.20407366666667, 159.612221, 0.9416818042520982, 0.9416818042520981, 3, 1, 4.452417, 13.357251, 0.01313274635953239, 0.01313274635953239, 236.886722062719, 2131.980498564471, 0.01236686828662963, 0.01236686828662963),
    @("ECs", "Mmrn2", "Cd93", "Resolving-Mac", 3, 1, 53.20407366666667, 159.612221, 0.9416818042520982, 0.9416818042520981, 3, 1, 53.03808999999999, 159.11427, 0.156439925407717, 0.156439925407717, 2821.842447499296, 25396.58202749367, 0.1473166312150026, 0.1473166312150026),
    @("FAPs", "Mmrn2", "Cd93", "ECs", 3, 1, 1.641018666666667, 4.923056, 0.02904509584209167, 0.02904509584209167, 3, 1, 281.0920463333333, 843.2761389999999, 0.8291026083535286, 0.8291026083535286, 461.2772950845315, 4151.495655760784, 0.02408136472255644, 0.02408136472255644),
    @("FAPs", "Mmrn2", "Cd93", "FAPs", 3, 1, 1.641018666666667, 4.923056, 0.02904509584209167, 0.02904509584209167, 2, 0.6666666666666666, 0.449122, 1.347366, 0.001324719879221983, 0.001324719879221983, 0.7370175856106667, 6.633158270496, [double]"3.847661585592661E-05", [double]"3.847661585592661E-05"),
    @("FAPs", "Mmrn2", "Cd93", "MuSCs", 3, 1, 1.641018666666667, 4.923056, 0.02904509584209167, 0.02904509584209167, 3, 1, 4.452417, 13.357251, 0.01313274635953239, 0.01313274635953239, 7.306499408783999, 65.758494679056, 0.0003814418766824987, 0.0003814418766824988),
    @("FAPs", "Mmrn2", "Cd93", "Resolving-Mac", 3, 1, 1.641018666666667, 4.923056, 0.02904509584209167, 0.02904509584209167, 3, 1, 53.03808999999999, 159.11427, 0.156439925407717, 0.156439925407717, 87.03649573434664, 783.3284616091198, 0.004543812626996812, 0.004543812626996812),
    @("MuSCs", "Mmrn2", "Cd93", "ECs", 3, 1, 1.653900666666667, 4.961702, 0.02927309990581012, 0.02927309990581012, 3, 1, 281.0920463333333, 843.2761389999999, 0.8291026083535286, 0.8291026083535286, 464.8983228253975, 4184.084905428577, 0.0242704034865006, 0.0242704034865006),
    @("MuSCs", "Mmrn2", "Cd93", "FAPs", 3, 1, 1.653900666666667, 4.961702, 0.02927309990581012, 0.02927309990581012, 2, 0.6666666666666666, 0.449122, 1.347366, 0.001324719879221983, 0.001324719879221983, 0.7428031752146667, 6.685228576932, [double]"3.877865737167784E-05", [double]"3.877865737167784E-05"),
    @("MuSCs", "Mmrn2", "Cd93", "MuSCs", 3, 1, 1.653900666666667, 4.961702, 0.02927309990581012, 0.02927309990581012, 3, 1, 4.452417, 13.357251, 0.01313274635953239, 0.01313274635953239, 7.363855444577999, 66.27469900120199, 0.0003844361962202557, 0.0003844361962202557),
    @("MuSCs", "Mmrn2", "Cd93", "Resolving-Mac", 3, 1, 1.653900666666667, 4.961702, 0.02927309990581012, 0.02927309990581012, 3, 1, 53.03808999999999, 159.11427, 0.156439925407717, 0.156439925407717, 87.71973240972665, 789.4775916875399, 0.004579481565717582, 0.004579481565717583)
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowNum = $i + 2
    $rowData = $newValues[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowData[$j]
    }
}
